# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# Updates the recalculated calibration outputs (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN)
# for the affected rows in the risk-neutral extrapolation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @{ D = 116517.3611390926;  E = 0.0225526219031326;   F = 0.1619315680856225;  G = -0.9980679205896916; H = 10.27151143191923 }
    5  = @{ D = 118190.1850686447;  E = 0.004616184557523723; F = 0.2045734683473742;  G = -0.7471488745757884; H = 8.140619230617689 }
    6  = @{ D = 118589.4665493969;  E = -0.01073324829956123; F = 0.2461963849384141;  G = -1.362396611284616;  H = 11.44906056822656 }
    8  = @{ D = 120272.6170416691;  E = -0.03189599974531242; F = 0.2120320862083005;  G = -0.9301164525105129; H = 6.919902804251932 }
    9  = @{ D = 121836.19929144;    E = -0.06405189092210928; F = 0.3525291549876354;  G = -2.114120310140174;  H = 13.76292796154375 }
    10 = @{ D = 123313.7596198474; E = -0.1051361967047122;   F = 0.4453149052701786;  G = -1.929441587427766;  H = 9.80963810794877 }
    13 = @{ D = 115820.8310141887; E = 0.09635463804703032;   F = 0.1123356836587552;  G = -0.7277257001707244; H = 6.555853308919397 }
    17 = @{ D = 115823.0423619805; E = 0.1766734332596782;    F = 0.09935124558109336; G = -0.5818044104059169; H = 6.353284726074335 }
    20 = @{ D = 116694.7614515533; E = 0.05247086805839964;   F = 0.134702381573318;   G = -0.5226675737249257; H = 6.782765781999012 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
